{"js": "// Apply the Klondike Fever review edits via the Word JavaScript API.\n// Each entry is an exact, unique (or intentionally repeated) substring\n// replacement performed with Body.search() + Range.insertText(text, \"Replace\").\n\nconst replacements = [\n  {\n    from: \"Play Klondike Fever Free: Slot Game Review\",\n    to: \"Play Klondike Fever Free: Engaging Mining-themed Slot Game\",\n  },\n  {\n    from: \"Engaging gameplay with various special features\",\n    to: \"Variety of special features (Wilds, Scatters, Magic Spin, Mine Bonus)\",\n  },\n  {\n    from: \"Tasteful and immersive graphics\",\n    to: \"Engaging mining-themed setting and graphics\",\n  },\n  {\n    from: \"Easy-to-understand game mechanics\",\n    to: \"25 winning lines and 5x3 game grid\",\n  },\n  {\n    from: \"25 winning lines with RTP of 95%\",\n    to: \"Straightforward game mechanics\",\n  },\n  {\n    from: \"No background music\",\n    to: \"No background music, only mining sounds\",\n  },\n  {\n    from: \"Magic Spin feature only available with special symbols\",\n    to: \"RTP of 95% is average\",\n  },\n  {\n    from:\n      \"Discover Klondike Fever, an engaging slot game with special features and immersive graphics. Play free now and seek fortunes in the mining world.\",\n    to: \"Experience the excitement of Klondike Fever with its special features. Play now for free!\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { from, to } of replacements) {\n  const found = body.search(from, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Klondike Fever review edits via Word COM interop (Find/Replace).\n# wdReplaceAll = 2 handles both occurrences of the title text in one call;\n# it is also safe for the other, already-unique strings.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"Play Klondike Fever Free: Slot Game Review\"; New = \"Play Klondike Fever Free: Engaging Mining-themed Slot Game\" },\n    @{ Old = \"Engaging gameplay with various special features\"; New = \"Variety of special features (Wilds, Scatters, Magic Spin, Mine Bonus)\" },\n    @{ Old = \"Tasteful and immersive graphics\"; New = \"Engaging mining-themed setting and graphics\" },\n    @{ Old = \"Easy-to-understand game mechanics\"; New = \"25 winning lines and 5x3 game grid\" },\n    @{ Old = \"25 winning lines with RTP of 95%\"; New = \"Straightforward game mechanics\" },\n    @{ Old = \"No background music\"; New = \"No background music, only mining sounds\" },\n    @{ Old = \"Magic Spin feature only available with special symbols\"; New = \"RTP of 95% is average\" },\n    @{ Old = \"Discover Klondike Fever, an engaging slot game with special features and immersive graphics. Play free now and seek fortunes in the mining world.\"; New = \"Experience the excitement of Klondike Fever with its special features. Play now for free!\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$r.New, 2)\n}\n"}
